# Fixed description storage in JSON
# The sheet stores a set of "field name" / "field value (html)" pairs in
# columns A and B. A new "desc" field is inserted between "date" and
# "name", pushing "name" and "sity" down by one row, and the HTML that
# used to be split across the country/info cells is re-split so that the
# country info and the (new) description info each get their own cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, shift the existing "name" (row 5) and "sity" (row 6) rows down
# by one row to make room for the new "desc" row, copying both value and
# formatting so the A7 label cell keeps the same style as the rest of the
# column A header cells.
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("B6").Copy($ws.Range("B7"))

$ws.Range("A5").Copy($ws.Range("A6"))
$ws.Range("B5").Copy($ws.Range("B6"))

# Update the top "One..." id string.
$ws.Range("B1").Value = "One2684897887"

# Column B html values keep the same address text.
$ws.Range("B2").Value = "<strong>Адрес: </strong>Adderss<br>"

# The country info used to be bundled together with the "info" dropdown
# html in one shared string; now it lives on its own in B3.
$ws.Range("B3").Value = "<strong> Страна: </strong>Россия<br>"

# Date value changed from 24-04-2020 to 21-04-2020.
$ws.Range("B4").Value = "<strong> Дата: 21-04-2020</strong></p>"

# New "desc" row: label in A5, and the (fixed) description dropdown html
# in B5 -- "dfhxsn" was replaced with "description", and the trailing
# country <strong> fragment was removed since it now lives in B3.
$ws.Range("A5").Value = "desc"
$ws.Range("B5").Value = '<p> <li style="list-style-type: none;" >                      <a href="#" class="" style="padding: 0px" data-toggle="dropdown" role="button"                      aria-haspopup="true" aria-expanded="false"><strong>Инфо</strong>                      <span class="caret"></span></a> <ul class="dropdown-menu">                   <li>description</li></ul> </li>'

# Rows 6 ("name") and 7 ("sity") keep their original B values, already
# copied above -- just make sure the A labels are correct.
$ws.Range("A6").Value = "name"
$ws.Range("A7").Value = "sity"
